# Quarterly financials update: insert 3 new columns (new quarter + two more recent
# quarters) before column D, shifting the existing quarterly columns D:K to G:N,
# and populate the 3 new columns with the latest reported figures. Also corrects a
# couple of historical "Capital Expenditures" cash-flow figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert three blank columns before column D. This automatically shifts the
#    existing D:K data (and its formatting) right to G:N.
# ---------------------------------------------------------------------------
$ws.Range("D:F").Insert()

# ---------------------------------------------------------------------------
# 2. Apply number formats to the new D:F columns.
#    Rows 7, 38 and 80 hold the period-ending dates; everything else is a
#    plain thousands-formatted number.
# ---------------------------------------------------------------------------
$dateFormat = "[$-409]d\-mmm\-yy;@"
$numberFormat = "#,##0"

$dateRows = @(7, 38, 80)
foreach ($r in $dateRows) {
    $ws.Range("D" + $r + ":F" + $r).NumberFormat = $dateFormat
}

$numberRows = @(8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,
                39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,
                81,82,83,84,85,86,87,88,89,90,91,92,93,94,95,96,97,98,99,100,101,102)
foreach ($r in $numberRows) {
    $ws.Range("D" + $r + ":F" + $r).NumberFormat = $numberFormat
}

# ---------------------------------------------------------------------------
# 3. New quarterly data for columns D (latest quarter), E and F.
#    Key: row number -> array(D, E, F). "NA" denotes the literal text used
#    elsewhere in the sheet for not-available figures; $null leaves the cell
#    blank (used for section header rows).
# ---------------------------------------------------------------------------
$newData = @{
    7  = @(43465, 43373, 43281)
    8  = @(115600, 117800, 113500)
    9  = @(28500, 28300, 29900)
    10 = @(87100, 89500, 83600)
    11 = @($null, $null, $null)
    12 = @("NA", "NA", "NA")
    13 = @(0, 0, 0)
    14 = @(210700, -116400, "NA")
    15 = @(27000, 27000, 26700)
    16 = @($null, $null, $null)
    17 = @(276500, -51000, 64800)
    18 = @(-160900, 168800, 48700)
    19 = @($null, $null, $null)
    20 = @(-700, -21000, -19800)
    21 = @(-134700, 174700, 55600)
    22 = @(19300, 20500, 23000)
    23 = @(-181000, 127200, 5800)
    24 = @(0, 0, 0)
    25 = @(0, 0, 0)
    26 = @(-181000, 127200, 5800)
    27 = @(-181000, 127200, 5800)
    28 = @(0, 0, 0)
    29 = @(0, 0, 0)
    30 = @(0, 0, 0)
    31 = @(0, 0, 0)
    32 = @(700, 21000, 19800)
    33 = @(-181000, 127200, 5800)
    34 = @(0, 0, 0)
    35 = @(-181000, 127200, 5800)

    38 = @(43465, 43373, 43281)
    39 = @($null, $null, $null)
    40 = @($null, $null, $null)
    41 = @(77300, 80100, 72400)
    42 = @(0, 0, 0)
    43 = @(27200, 21900, 50700)
    44 = @(8900, 9700, 9100)
    45 = @(6400, 7900, 10500)
    46 = @(119800, 119700, 142700)
    47 = @(36300, 48000, 45300)
    48 = @(2480300, 2717300, 2743900)
    49 = @(0, 0, 0)
    50 = @(0, 0, 0)
    51 = @(0, 0, 0)
    52 = @(43500, 36400, 34900)
    53 = @(0, 0, 0)
    54 = @(2679800, 2921300, 2966800)
    55 = @($null, $null, $null)
    56 = @($null, $null, $null)
    57 = @(46300, 51300, 16300)
    58 = @(113800, 118500, 134900)
    59 = @(62700, 69600, 38400)
    60 = @(222700, 239400, 189600)
    61 = @(1508100, 1527400, 2153000)
    62 = @(258200, 272700, 50300)
    63 = @(0, 0, 0)
    64 = @(0, 0, 0)
    65 = @(0, 0, 0)
    66 = @(1989000, 2039500, 2392900)
    67 = @($null, $null, $null)
    68 = @(0, 0, 0)
    69 = @(0, 0, 0)
    70 = @(0, 0, 0)
    71 = @(0, 0, 0)
    72 = @(81800, 262800, 135600)
    73 = @(0, 0, 0)
    74 = @(0, 0, 0)
    75 = @(0, 0, 0)
    76 = @(690900, 881800, 573900)
    77 = @(0, 0, 0)

    80 = @(43465, 43373, 43281)
    81 = @(-181000, 127200, 5800)
    82 = @($null, $null, $null)
    83 = @(27000, 27000, 26700)
    84 = @(0, 0, 0)
    85 = @(0, 0, 0)
    86 = @(0, 0, 0)
    87 = @(0, 0, 0)
    88 = @(0, 0, 0)
    89 = @(53800, 55300, 20000)
    90 = @($null, $null, $null)
    91 = @(-3300, -400, -1000)
    92 = @(0, 0, 0)
    93 = @(0, 0, 0)
    94 = @(-6200, -400, -1000)
    95 = @($null, $null, $null)
    96 = @(0, 0, 0)
    97 = @(0, 0, 0)
    98 = @(0, 0, 0)
    99 = @(0, 0, 0)
    100 = @(-50500, -50000, -6800)
    101 = @(0, 0, 0)
    102 = @(-2900, 4900, 12300)
}

foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    $ws.Cells.Item($r, 4).Value = $vals[0]
    $ws.Cells.Item($r, 5).Value = $vals[1]
    $ws.Cells.Item($r, 6).Value = $vals[2]
}

# ---------------------------------------------------------------------------
# 4. A small correction to the "Capital Expenditures" row: the quarter that
#    shifted into column G (and the two quarters after it) were updated from
#    0 to their corrected reported values.
# ---------------------------------------------------------------------------
$ws.Cells.Item(91, 7).Value = -700
$ws.Cells.Item(91, 8).Value = "NA"
$ws.Cells.Item(91, 9).Value = "NA"
$ws.Cells.Item(91, 10).Value = "NA"
